# 1st iteration (#141)
# Update the "Experimental" flag to the literal text "true" and bump the
# "Date" metadata timestamp on the Metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7: Property = "Experimental", Value = "true"
# Writing "true" straight to Range.Value would be auto-coerced to a Boolean
# by Excel's type inference, but the source data needs it stored as text.
# Build the text in a scratch cell via a formula (so no quote-prefix style
# gets attached), copy/paste its *value* into B7, then clean up the scratch
# cell so it leaves no trace in the sheet.
$ws.Range("D1").Formula = "=""true"""
$ws.Range("D1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("D1").Clear()

# Row 8: Property = "Date", Value = updated ISO timestamp
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
